$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 141-21 84th Drive in Jamaica, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/10/permits-filed-for-141-21-84th-drive-in-jamaica-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a four-story residential building at 141-21 84th Drive in <a href="https://newyorkyimby.com/neighborhoods/jamaica">Jamaica</a>, Queens. Located between 85th Road and Burden Crescent, the lot is near the Briarwood subway station, served by the E and F trains. Jacob Ashkenazie is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2025-10-04T10:30:36+00:00"
$ws.Range("E2").Value = "Sat, 04 Oct 2025 10:30:36 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Jamaica"
$ws.Range("H2").Value = ""
